# Generate Report for Handoff
# Update the localization status report with the newly generated file's
# GUID-based name, content hash and refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$newGuid = "0804cc9e-3e7a-46d5-ad3e-5ab4dd460376"
$newHash = "22bdf1ce61407596ce5696f2971039d13c50dbd5"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-31 08:03:18"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-31 08:03:00"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-31 08:03:18"

